# Manual update + fixes to the scraped-data sheet:
#  - daily "nothing relevant" rows (2025-02-18 .. 2025-02-25) are collapsed to a
#    single occurrence count of 1 each, and their dates are corrected
#  - 2025-02-21/02-24/02-25 "nothing relevant" rows are replaced by real
#    term/page rows for 2025-02-26, and every subsequent page-number (C) /
#    occurrence-count (D) value on 2025-02-26 is refreshed
#  - the trailing duplicate rows (42-45) are removed, shrinking A1:D45 to A1:D41
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure date-like strings in column A stay as text, not auto-converted to dates
$ws.Range("A2:A41").NumberFormat = "@"

$rows = @(
    @(2, '2025-02-18', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(3, '2025-02-19', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(4, '2025-02-20', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(5, '2025-02-20', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(6, '2025-02-21', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(7, '2025-02-24', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(8, '2025-02-25', 'Rien ne nous concerne aujourd''hui !', $null, 1),
    @(9, '2025-02-26', 'autorisation', 1, 1),
    @(10, '2025-02-26', 'service', 1, 5),
    @(11, '2025-02-26', 'service', 2, 4),
    @(12, '2025-02-26', 'autorisation', 2, 1),
    @(13, '2025-02-26', 'gouvernement', 2, 1),
    @(14, '2025-02-26', 'service', 3, 6),
    @(15, '2025-02-26', 'gouvernement', 3, 1),
    @(16, '2025-02-26', 'service', 4, 11),
    @(17, '2025-02-26', 'autorisation', 5, 3),
    @(18, '2025-02-26', 'service', 6, 6),
    @(19, '2025-02-26', 'service', 7, 5),
    @(20, '2025-02-26', 'service', 8, 3),
    @(21, '2025-02-26', 'service', 9, 3),
    @(22, '2025-02-26', 'autorisation', 9, 3),
    @(23, '2025-02-26', 'gouvernement', 13, 3),
    @(24, '2025-02-26', 'service', 13, 2),
    @(25, '2025-02-26', 'gouvernement', 14, 7),
    @(26, '2025-02-26', 'service', 14, 1),
    @(27, '2025-02-26', 'autorisation', 14, 1),
    @(28, '2025-02-26', 'service', 19, 2),
    @(29, '2025-02-26', 'service', 20, 7),
    @(30, '2025-02-26', 'service', 21, 2),
    @(31, '2025-02-26', 'service', 22, 4),
    @(32, '2025-02-26', 'service', 23, 2),
    @(33, '2025-02-26', 'service', 25, 3),
    @(34, '2025-02-26', 'service', 26, 3),
    @(35, '2025-02-26', 'service', 27, 2),
    @(36, '2025-02-26', 'service', 29, 4),
    @(37, '2025-02-26', 'service', 30, 2),
    @(38, '2025-02-26', 'service', 31, 3),
    @(39, '2025-02-26', 'service', 32, 4),
    @(40, '2025-02-26', 'service', 33, 1),
    @(41, '2025-02-26', 'service', 34, 2)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    if ($null -eq $row[3]) {
        $ws.Cells.Item($r, 3).Value = ""
    } else {
        $ws.Cells.Item($r, 3).Value = [double]$row[3]
    }
    $ws.Cells.Item($r, 4).Value = [double]$row[4]
}

# Remove the now-obsolete trailing rows (42-45), shrinking the used range to A1:D41
$ws.Range("A42:D45").EntireRow.Delete()